$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.685592333333333
$ws.Range("H2").Value = 17.056777
$ws.Range("I2").Value = 0.1188473284691575
$ws.Range("J2").Value = 0.1188473284691575
$ws.Range("M2").Value = 3.759736666666667
$ws.Range("N2").Value = 11.27921
$ws.Range("O2").Value = 0.0683751702595819
$ws.Range("P2").Value = 0.06837517025958188
$ws.Range("Q2").Value = 21.37632996735222
$ws.Range("R2").Value = 192.38696970617
$ws.Range("S2").Value = 0.008126206318975102
$ws.Range("T2").Value = 0.0081262063189751
$ws.Range("G3").Value = 5.685592333333333
$ws.Range("H3").Value = 17.056777
$ws.Range("I3").Value = 0.1188473284691575
$ws.Range("J3").Value = 0.1188473284691575
$ws.Range("O3").Value = 0.6514180024294648
$ws.Range("P3").Value = 0.6514180024294647
$ws.Range("Q3").Value = 203.6547201819112
$ws.Range("R3").Value = 1832.892481637201
$ws.Range("S3").Value = 0.07741928930545708
$ws.Range("T3").Value = 0.07741928930545706
$ws.Range("G4").Value = 5.685592333333333
$ws.Range("H4").Value = 17.056777
$ws.Range("I4").Value = 0.1188473284691575
$ws.Range("J4").Value = 0.1188473284691575
$ws.Range("O4").Value = 0.2802068273109533
$ws.Range("P4").Value = 0.2802068273109533
$ws.Range("Q4").Value = 87.60188204232556
$ws.Range("R4").Value = 788.41693838093
$ws.Range("S4").Value = 0.03330183284472538
$ws.Range("T4").Value = 0.03330183284472537
$ws.Range("I5").Value = 0.622926875404983
$ws.Range("J5").Value = 0.6229268754049829
$ws.Range("M5").Value = 3.759736666666667
$ws.Range("N5").Value = 11.27921
$ws.Range("O5").Value = 0.0683751702595819
$ws.Range("P5").Value = 0.06837517025958188
$ws.Range("Q5").Value = 112.0419836584233
$ws.Range("R5").Value = 1008.37785292581
$ws.Range("S5").Value = 0.04259273116508507
$ws.Range("T5").Value = 0.04259273116508506
$ws.Range("I6").Value = 0.622926875404983
$ws.Range("J6").Value = 0.6229268754049829
$ws.Range("O6").Value = 0.6514180024294648
$ws.Range("P6").Value = 0.6514180024294647
$ws.Range("S6").Value = 0.4057857808359421
$ws.Range("T6").Value = 0.405785780835942
$ws.Range("I7").Value = 0.622926875404983
$ws.Range("J7").Value = 0.6229268754049829
$ws.Range("O7").Value = 0.2802068273109533
$ws.Range("P7").Value = 0.2802068273109533
$ws.Range("S7").Value = 0.1745483634039558
$ws.Range("T7").Value = 0.1745483634039557
$ws.Range("I8").Value = 0.2582257961258595
$ws.Range("J8").Value = 0.2582257961258594
$ws.Range("M8").Value = 3.759736666666667
$ws.Range("N8").Value = 11.27921
$ws.Range("O8").Value = 0.0683751702595819
$ws.Range("P8").Value = 0.06837517025958188
$ws.Range("Q8").Value = 46.44546827572222
$ws.Range("R8").Value = 418.0092144815
$ws.Range("S8").Value = 0.01765623277552173
$ws.Range("T8").Value = 0.01765623277552172
$ws.Range("I9").Value = 0.2582257961258595
$ws.Range("J9").Value = 0.2582257961258594
$ws.Range("O9").Value = 0.6514180024294648
$ws.Range("P9").Value = 0.6514180024294647
$ws.Range("S9").Value = 0.1682129322880656
$ws.Range("T9").Value = 0.1682129322880656
$ws.Range("I10").Value = 0.2582257961258595
$ws.Range("J10").Value = 0.2582257961258594
$ws.Range("O10").Value = 0.2802068273109533
$ws.Range("P10").Value = 0.2802068273109533
$ws.Range("S10").Value = 0.07235663106227215
$ws.Range("T10").Value = 0.07235663106227212
